$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.962.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.649.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.01'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0756'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.973.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '193.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.46%  '
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.52'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.895'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.121.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.532'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.795'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '55.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0937'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("E51").Value = '  +0.36%  '
